$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = "새로운 루다를 지탱하는 모델 서빙 아키텍처 — 1편: A/B 테스트를 위한 구조 설계"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/serving-architecture-1/"

$ws.Range("D32").Value = "[Airflow] Xcom을 이용한 task간 변수 전달"
$ws.Range("E32").Value = "https://dodonam.tistory.com/402"

$ws.Range("D36").Value = "Transformer-based Anomaly Detection in Multivariate Time Series"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/394"

$ws.Range("D51").Value = "[python] 비공개 속성, 게터, 세터, @property"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EB%B9%84%EA%B3%B5%EA%B0%9C-%EC%86%8D%EC%84%B1-%EA%B2%8C%ED%84%B0-%EC%84%B8%ED%84%B0-property"
